# The 2011 edit: the header sub-row (units: (m3/s), (MW), (GWh)) that sat
# between the column-label row and the first data row is removed. All the
# data rows shift up by one; Excel drops the now-unused "(MW)"/"(GWh)"
# shared strings and recomputes the used range (dimension) automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the units sub-header row (row 2); everything below shifts up.
$ws.Rows.Item(2).Delete()

# Reflect the new active selection on the first data row, like in the
# authored workbook (A2:K2 instead of the old A1:K1).
$ws.Range("A2:K2").Select()
